$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - "Expandable subitem clic go to test" (week 14)
$ws.Range("M21").Value = 4
$ws.Range("AA21").Value = "Realizado el 17/10/18"
$ws.Range("AB21").Value = 20

# Row 22 - "Back up's Git" (week 15)
$ws.Range("M22").Value = 0.5

# Row 23 - "Documentación código fuente" (week 16)
$ws.Range("J23").Value = 16
$ws.Range("M23").Value = 12
$ws.Range("AB23").Value = 100

# Row 24 - new task (week 17)
$ws.Range("C24").Value = "SharedPreferences para el JSON desde el ExpandableListView"
$ws.Range("M24").Value = 6
$ws.Range("AB24").Value = 10

# Row 25 - new task (week 18)
$ws.Range("C25").Value = "Create function: check if exit json. If not exit show message"
$ws.Range("M25").Value = 6
$ws.Range("AB25").Value = 6

# Update active cell selection to reflect where the author ended up editing
$ws.Range("M22").Select()
